# Applies the edit described by the diff:
# 1. Splits the run containing "he is sitting in council" with gramStart/gramEnd proofErr marks.
# 2. Splits the "These must, indeed..." paragraph similarly (two proofErr pairs).
# 3. Splits the "So the two pretended weavers..." paragraph similarly (two proofErr pairs).
# 4. Rewrites the final (bookmark-only) paragraph to add a divider and new story text
#    around the existing _GoBack bookmark, and removes the stray w:hint="cs"/w:cs from pPr/rPr.
# 5. Appends a brand-new paragraph with more story text after it.

$d = $word.ActiveDocument

# --- Paragraph 2: "Many years ago, ... wardrobe." ---
$para2Xml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>Many years ago, there was an Emperor, who was so excessively fond of new clothes, that he</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>spent all his money in dress. He did not trouble himself in the least about his soldiers; nor did</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>he care to go either to the theatre or the chase, except for the opportunities then afforded him</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>for displaying his new clothes. He had a different suit for each hour of the day; and as of any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>other king or emperor, one is accustomed to say, "</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>he</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> is sitting in council," it was always said of</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>him, "The Emperor is sitting in his wardrobe."</w:t></w:r></w:p>
'@
$d.Paragraphs(2).Range.InsertXML($para2Xml)

# --- Paragraph 5: "These must, indeed, ... directly." ---
$paraTheseXml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">--------------------------------------------------------------------------------------------------------------- "These must, indeed, be splendid clothes!" thought the Emperor. "Had I such a suit, I might </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>at  once</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> find out what men in my realms are unfit for their office, and also be able to distinguish  the wise from the foolish! This stuff must be woven for me immediately." And he caused </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>large  sums</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> of money to be given to both the weavers in order that they might begin their work  directly. </w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($paraTheseXml)

# --- Paragraph 6: "So the two pretended weavers ... late at night." ---
$paraWeaversXml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">So the two pretended weavers set up two looms, and affected to work very busily, though </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>in  reality</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> they did nothing at all. They asked for the most delicate silk and the purest gold thread</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>;  put</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> both into their own knapsacks; and then continued their pretended work at the empty  looms until late at night. </w:t></w:r></w:p>
'@
$d.Paragraphs(6).Range.InsertXML($paraWeaversXml)

# --- Paragraph 7: previously-empty bookmark paragraph -> dashes + bookmark + new text,
#     plus a brand new paragraph 8 appended after it. Both supplied in one InsertXML call
#     so the trailing sectPr / document end is not disturbed. ---
$para7Xml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>-----------------------------------------------------------------------------------------------</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>---------------</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> "I should like to know how the weavers are getting on with my cloth," said the Emperor to  himself, after some little time had elapsed; he was, however, rather embarrassed, when he  remembered that a simpleton, or one unfit for his office, would be unable to see the  manufacture. To be sure, he thought he had nothing to risk in his own person; but yet, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>he  would</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> prefer sending somebody else, to bring him intelligence about the weavers, and their  work, before he troubled himself in the affair. All the people throughout the city had heard </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>of  the</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> wonderful property the cloth was to possess; and all were anxious to learn how wise, or  how ignorant, their neighbors might prove to be. </w:t></w:r></w:p>
'@
$para8Xml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:cs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">"I will send my faithful old minister to the weavers," said the Emperor at last, after </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>some  deliberation</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>, "he will be best able to see how the cloth looks; for he is a man of sense, and no  one can be more suitable for his office than he is."</w:t></w:r></w:p>
'@
$d.Paragraphs(7).Range.InsertXML($para7Xml + $para8Xml)
